# Update "想去人数" (want-to-go count) figures and one cover image URL
# across the workbook's sheets, per the upstream data refresh commit
# "Update gh-pages to output generated at 456a3b4".
#
# Sheet 1 = 展览 (Exhibitions)
# Sheet 2 = 演出 (Performances)
# Sheet 3 = 本地生活 (Local Life)
# Sheet 4 = 全部类型 (All Types - combined view of the above three sheets)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 63
$ws1.Range("F7").Value = 974
$ws1.Range("F8").Value = 937
$ws1.Range("F13").Value = 920
$ws1.Range("F15").Value = 3902
$ws1.Range("F16").Value = 1155
$ws1.Range("F18").Value = 2585
$ws1.Range("F20").Value = 1077
$ws1.Range("F21").Value = 3566
$ws1.Range("F22").Value = 753
$ws1.Range("F23").Value = 838
$ws1.Range("F26").Value = 110
$ws1.Range("F27").Value = 832
$ws1.Range("F28").Value = 166
$ws1.Range("F29").Value = 441
$ws1.Range("F30").Value = 198
$ws1.Range("F32").Value = 1330
$ws1.Range("F33").Value = 1944
$ws1.Range("F35").Value = 37
$ws1.Range("F37").Value = 586
$ws1.Range("F39").Value = 15
$ws1.Range("I39").Value = "//i0.hdslb.com/bfs/openplatform/202404/J7M4PT141713239020893.jpeg"
$ws1.Range("F41").Value = 232

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F11").Value = 19

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 425

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 425
$ws4.Range("F5").Value = 63
$ws4.Range("F6").Value = 974
$ws4.Range("F7").Value = 937
$ws4.Range("F14").Value = 920
$ws4.Range("F16").Value = 3902
$ws4.Range("F17").Value = 1155
$ws4.Range("F20").Value = 2585
$ws4.Range("F21").Value = 1077
$ws4.Range("F22").Value = 3566
$ws4.Range("F23").Value = 753
$ws4.Range("F24").Value = 838
$ws4.Range("F30").Value = 19
$ws4.Range("F31").Value = 110
$ws4.Range("F33").Value = 832
$ws4.Range("F34").Value = 166
$ws4.Range("F35").Value = 441
$ws4.Range("F36").Value = 198
$ws4.Range("F38").Value = 1330
$ws4.Range("F39").Value = 1944
$ws4.Range("F43").Value = 37
$ws4.Range("F44").Value = 586
$ws4.Range("F46").Value = 15
$ws4.Range("I46").Value = "//i0.hdslb.com/bfs/openplatform/202404/J7M4PT141713239020893.jpeg"
$ws4.Range("F48").Value = 232
